$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 36
$ws.Range("A36").Value = 112073630
$ws.Range("B36").Value = 89072
$ws.Range("D36").Value = 'LC'
$ws.Range("E36").Value = 256703
$ws.Range("F36").Value = 'Tallfingersvamp'
$ws.Range("G36").Value = 'Ramaria eosanguinea'
$ws.Range("H36").Value = 'R.H.Petersen'
$ws.Range("I36").Value = '''1'
$ws.Range("J36").Value = 'fruktkroppar'
$ws.Range("P36").Value = 'N om Nedre Tetvasseltjärnen, Dlr'
$ws.Range("Q36").Value = 491917
$ws.Range("R36").Value = 6785497
$ws.Range("Z36").Value = ""
$ws.Range("AB36").Value = ""
$ws.Range("AH36").Value = 'Sandtallskog'
$ws.Range("AJ36").Value = 'tall'
$ws.Range("AK36").Value = 'Pinus sylvestris'
$ws.Range("AO36").Value = 'Pinus sylvestris'
$ws.Range("AW36").Value = 'Janolof Hermansson'
$ws.Range("AX36").Value = 'Janolof Hermansson, Bengt Oldhammer, Bo karlstens, Birgitta Kvist'

# Row 37
$ws.Range("A37").Value = 112073422
$ws.Range("B37").Value = 89097
$ws.Range("D37").Value = 'LC'
$ws.Range("E37").Value = 233195
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = 'Ramaria neoformosa'
$ws.Range("H37").Value = 'sensu Schild'
$ws.Range("I37").Value = '''3'
$ws.Range("P37").Value = 'N om Nedre Tetvasseltjärnen, Dlr'
$ws.Range("Q37").Value = 491996
$ws.Range("R37").Value = 6785531
$ws.Range("S37").Value = 5
$ws.Range("Z37").Value = ""
$ws.Range("AB37").Value = ""
$ws.Range("AH37").Value = 'Sandtallskog'
$ws.Range("AJ37").Value = 'tall'
$ws.Range("AK37").Value = 'Pinus sylvestris'
$ws.Range("AO37").Value = 'Pinus sylvestris'
$ws.Range("AW37").Value = 'Janolof Hermansson'
$ws.Range("AX37").Value = 'Janolof Hermansson, Bengt Oldhammer, Bo karlstens, Birgitta Kvist'

# Row 38
$ws.Range("A38").Value = 112045302
$ws.Range("B38").Value = 90806
$ws.Range("E38").Value = 4361
$ws.Range("F38").Value = 'Orange taggsvamp'
$ws.Range("G38").Value = 'Hydnellum aurantiacum'
$ws.Range("H38").Value = '(Batsch:Fr.) P.Karst.'
$ws.Range("S38").Value = 5
$ws.Range("Z38").Value = '12:27'
$ws.Range("AB38").Value = '12:27'

# Row 39
$ws.Range("A39").Value = 112045414
$ws.Range("B39").Value = 90837
$ws.Range("D39").Value = 'NT'
$ws.Range("E39").Value = 5966
$ws.Range("F39").Value = 'Motaggsvamp'
$ws.Range("G39").Value = 'Sarcodon squamosus'
$ws.Range("H39").Value = '(Schaeff.) Quél.'
$ws.Range("I39").Value = ""
$ws.Range("P39").Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Range("Q39").Value = 492044
$ws.Range("R39").Value = 6785564
$ws.Range("S39").Value = 15
$ws.Range("Z39").Value = '12:36'
$ws.Range("AB39").Value = '12:36'
$ws.Range("AH39").Value = ""
$ws.Range("AJ39").Value = ""
$ws.Range("AK39").Value = ""
$ws.Range("AO39").Value = ""
$ws.Range("AW39").Value = 'Bo karlstens'
$ws.Range("AX39").Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson, Birgitta Kvist'

# Row 40
$ws.Range("A40").Value = 112045406
$ws.Range("B40").Value = 90830
$ws.Range("D40").Value = 'NT'
$ws.Range("E40").Value = 2059
$ws.Range("F40").Value = 'Skrovlig taggsvamp'
$ws.Range("G40").Value = 'Hydnellum scabrosum'
$ws.Range("H40").Value = '(Fr.) E.Larss., K.H.Larss. & Kõljalg'
$ws.Range("I40").Value = ""
$ws.Range("J40").Value = ""
$ws.Range("P40").Value = 'Nedre Tetvasseltjärnen (Nedre Tetvasseltjärnen), Dlr'
$ws.Range("Q40").Value = 492044
$ws.Range("R40").Value = 6785564
$ws.Range("S40").Value = 15
$ws.Range("Z40").Value = '12:36'
$ws.Range("AB40").Value = '12:36'
$ws.Range("AH40").Value = ""
$ws.Range("AJ40").Value = ""
$ws.Range("AK40").Value = ""
$ws.Range("AO40").Value = ""
$ws.Range("AW40").Value = 'Bo karlstens'
$ws.Range("AX40").Value = 'Bo karlstens, Bengt Oldhammer, Janolof Hermansson, Birgitta Kvist'
